$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-21"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 05-21)"

# Update the June (row 6) figure for the "through" column (I)
$ws.Range("I6").Value = 76

# Update the Total (row 14) figure for the "through" column (I)
$ws.Range("I14").Value = 628
